# Update the cryptos price/volume table with the latest scrape.
# Column D values that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the sheet's original
# inline-string formatting, e.g. "1.005" instead of numeric 1.005).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.474.51"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").Value = "1.722.23"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +1.03%  "

$ws.Range("D5").Value = "'225.07"
$ws.Range("E5").Value = "  -1.95%  "

$ws.Range("D6").Value = "'0.5334"
$ws.Range("E6").Value = "  -1.35%  "

$ws.Range("E7").Value = "  +1.00%  "

$ws.Range("D8").Value = "'0.2657"
$ws.Range("E8").Value = "  -2.02%  "

$ws.Range("D9").Value = "'0.06581"
$ws.Range("E9").Value = "  +1.47%  "

$ws.Range("D10").Value = "'21.31"
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").Value = "'0.07665"
$ws.Range("E11").Value = "  -0.89%  "

$ws.Range("D12").Value = "'4.603"
$ws.Range("E12").Value = "  -2.53%  "

$ws.Range("D13").Value = "1.731.76"
$ws.Range("E13").Value = "  +2.72%  "

$ws.Range("D14").Value = "1.960.17"
$ws.Range("E14").Value = "  +2.61%  "

$ws.Range("D15").Value = "'0.5787"
$ws.Range("E15").Value = "  -3.14%  "

$ws.Range("D16").Value = "0.0₅8261"
$ws.Range("E16").Value = "  -0.33%  "

$ws.Range("D17").Value = "'67.76"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").Value = "27.468.30"
$ws.Range("E18").Value = "  +2.30%  "

$ws.Range("D19").Value = "'217.75"
$ws.Range("E19").Value = "  +3.99%  "

$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("D21").Value = "'4.725"
$ws.Range("E21").Value = "  -1.31%  "

$ws.Range("E22").Value = "  -3.63%  "

$ws.Range("E23").Value = "  -2.63%  "

$ws.Range("E24").Value = "  +1.07%  "

$ws.Range("D25").Value = "'143.18"
$ws.Range("E25").Value = "  -2.37%  "

$ws.Range("D26").Value = "'1.738"
$ws.Range("E26").Value = "  +10.61%  "

$ws.Range("D27").Value = "'0.1228"
$ws.Range("E27").Value = "  -0.89%  "

$ws.Range("D28").Value = "'7.310"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("D29").Value = "'16.43"
$ws.Range("E29").Value = "  -0.76%  "

$ws.Range("D30").Value = "'0.05435"
$ws.Range("E30").Value = "  -4.59%  "

$ws.Range("E31").Value = "  -1.60%  "

$ws.Range("D32").Value = "'3.530"
$ws.Range("E32").Value = "  -2.38%  "

$ws.Range("D33").Value = "'3.431"
$ws.Range("E33").Value = "  -2.19%  "

$ws.Range("D34").Value = "'1.643"
$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("D35").Value = "'2.882"
$ws.Range("E35").Value = "  +2.02%  "

$ws.Range("D36").Value = "'0.9554"
$ws.Range("E36").Value = "  -2.71%  "

$ws.Range("D37").Value = "'2.429"
$ws.Range("E37").Value = "  +0.87%  "

$ws.Range("D38").Value = "'0.5901"
$ws.Range("E38").Value = "  +1.19%  "

$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("D40").Value = "'5.899"
$ws.Range("E40").Value = "  -1.27%  "

$ws.Range("D41").Value = "1.045.70"
$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").Value = "'0.8449"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("E43").Value = "  +0.99%  "

$ws.Range("D44").Value = "'101.24"
$ws.Range("E44").Value = "  -1.41%  "

$ws.Range("D45").Value = "1.866.59"
$ws.Range("E45").Value = "  +2.41%  "

$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  +9.68%  "

$ws.Range("D47").Value = "'58.49"
$ws.Range("E47").Value = "  -1.90%  "

$ws.Range("E48").Value = "  +4.34%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.160"
$ws.Range("E49").Value = "  +1.56%  "

$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").Value = "'1.003"
$ws.Range("E50").Value = "  +0.76%  "

$ws.Range("D51").Value = "'0.06601"
$ws.Range("E51").Value = "  +14.13%  "
